# 3/11/2016 -- St. Louis Competition -- Day 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update existing "Low Bar" entry to the close-range variant.
$ws.Range("B12").Value = "Low Bar (Close with /targetTrack)"

# Append the two new autonomous modes below the existing list.
$ws.Range("A13").Value = 21
$ws.Range("B13").Value = "Low Bar (Far with /targetTrack)"

$ws.Range("A14").Value = "default"
$ws.Range("B14").Value = "Corner Shot"

# Match styling of the rows above (same fill/font as rows 4-12).
$ws.Range("A13").Style = $ws.Range("A12").Style
$ws.Range("B13").Style = $ws.Range("B12").Style
$ws.Range("A14").Style = $ws.Range("A12").Style
$ws.Range("B14").Style = $ws.Range("B12").Style

# Update the selection to mirror the new used range / active cell.
$ws.Range("A1:B14").Select()
$ws.Range("B14").Activate()
